$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3335.75
$ws.Range("I32").Value = 2130.5
$ws.Range("J32").Value = 5746.25
$ws.Range("K32").Value = 2130.5
$ws.Range("L32").Value = 5746.25
$ws.Range("M32").Value = -1804.5
$ws.Range("N32").Value = -6398.25

$ws.Range("H69").Value = 22432
$ws.Range("J69").Value = 24073.834
$ws.Range("L69").Value = 72221.50199999999
$ws.Range("N69").Value = -73969.50199999999

$ws.Range("H72").Value = 22432
$ws.Range("J72").Value = 24073.834
$ws.Range("L72").Value = 216664.506
$ws.Range("N72").Value = -225400.506

$ws.Range("H132").Value = 4249.2104
$ws.Range("I132").Value = 4249.2104
$ws.Range("K132").Value = 12747.6312
$ws.Range("M132").Value = -10217.6312

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 25500
$ws.Range("J43").Value = 25500
$ws.Range("L43").Value = 25500
$ws.Range("N43").Value = -26126

$ws.Range("H61").Value = 6772.1577
$ws.Range("I61").Value = 5744.231
$ws.Range("K61").Value = 5744.231
$ws.Range("M61").Value = -5532.231

$ws.Range("H97").Value = 1720.9231
$ws.Range("I97").Value = 1071.091
$ws.Range("J97").Value = 5295
$ws.Range("K97").Value = 1071.091
$ws.Range("L97").Value = 5295
$ws.Range("M97").Value = -575.0909999999999
$ws.Range("N97").Value = -6287

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws.Range("H136").Value = 6772.1577
$ws.Range("I136").Value = 5744.231
$ws.Range("K136").Value = 17232.693
$ws.Range("M136").Value = -14682.693

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 937.5
$ws.Range("I11").Value = 933.3333
$ws.Range("J11").Value = 950
$ws.Range("K11").Value = 933.3333
$ws.Range("L11").Value = 950
$ws.Range("M11").Value = -793.3333
$ws.Range("N11").Value = -1230

$ws.Range("H86").Value = 7989.5654
$ws.Range("I86").Value = 3376.7856
$ws.Range("J86").Value = 15165
$ws.Range("K86").Value = 3376.7856
$ws.Range("L86").Value = 15165
$ws.Range("M86").Value = -2253.7856
$ws.Range("N86").Value = -17411

$ws.Range("H89").Value = 7989.5654
$ws.Range("I89").Value = 3376.7856
$ws.Range("J89").Value = 15165
$ws.Range("K89").Value = 16883.928
$ws.Range("L89").Value = 75825
$ws.Range("M89").Value = -11267.928
$ws.Range("N89").Value = -87057

$ws.Range("H134").Value = 10934.833
$ws.Range("I134").Value = 11339.8
$ws.Range("K134").Value = 34019.39999999999
$ws.Range("M134").Value = -31484.39999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 6100.2
$ws.Range("J62").Value = 6375.25
$ws.Range("L62").Value = 6375.25
$ws.Range("N62").Value = -7623.25

$ws.Range("H65").Value = 6100.2
$ws.Range("J65").Value = 6375.25
$ws.Range("L65").Value = 31876.25
$ws.Range("N65").Value = -38116.25

$ws.Range("H132").Value = 1478.25
$ws.Range("I132").Value = 1529.1666
$ws.Range("K132").Value = 4587.4998
$ws.Range("M132").Value = -2057.4998

$ws.Range("H134").Value = 3068.818
$ws.Range("I134").Value = 3480.5557
$ws.Range("K134").Value = 10441.6671
$ws.Range("M134").Value = -7906.667099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1547
$ws.Range("I5").Value = 1601
$ws.Range("K5").Value = 4803
$ws.Range("M5").Value = -4691

$ws.Range("H80").Value = 8659
$ws.Range("J80").Value = 8491
$ws.Range("L80").Value = 25473
$ws.Range("N80").Value = -27345

$ws.Range("H83").Value = 8659
$ws.Range("J83").Value = 8491
$ws.Range("L83").Value = 76419
$ws.Range("N83").Value = -85779

$ws.Range("H113").Value = 3805.3076
$ws.Range("I113").Value = 3096
$ws.Range("K113").Value = 9288
$ws.Range("M113").Value = -7118

$ws.Range("H132").Value = 2165.5
$ws.Range("I132").Value = 1831.3334
$ws.Range("J132").Value = 2499.6667
$ws.Range("K132").Value = 16482.0006
$ws.Range("L132").Value = 22497.0003
$ws.Range("M132").Value = -13952.0006
$ws.Range("N132").Value = -27557.0003

$ws.Range("H135").Value = 1547
$ws.Range("I135").Value = 1601
$ws.Range("K135").Value = 14409
$ws.Range("M135").Value = -11874

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3124
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 3124
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H97").Value = 820.73914
$ws.Range("I97").Value = 820.44446
$ws.Range("J97").Value = 821.8
$ws.Range("K97").Value = 820.44446
$ws.Range("L97").Value = 821.8
$ws.Range("M97").Value = -324.44446
$ws.Range("N97").Value = -1813.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1854.5714
$ws.Range("I22").Value = 836.6
$ws.Range("J22").Value = 4399.5
$ws.Range("K22").Value = 836.6
$ws.Range("L22").Value = 4399.5
$ws.Range("M22").Value = -541.6
$ws.Range("N22").Value = -4989.5

$ws.Range("H27").Value = 1854.5714
$ws.Range("I27").Value = 836.6
$ws.Range("J27").Value = 4399.5
$ws.Range("K27").Value = 836.6
$ws.Range("L27").Value = 4399.5
$ws.Range("M27").Value = -729.6
$ws.Range("N27").Value = -4613.5

$ws.Range("H70").Value = 33950
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 33950
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 33950
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -34490

$ws.Range("H73").Value = 33950
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 33950
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 33950
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -35822

$ws.Range("H122").Value = 3285.3333
$ws.Range("I122").Value = 3285.3333
$ws.Range("K122").Value = 9855.999899999999
$ws.Range("M122").Value = -7405.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 71333.336
$ws.Range("J63").Value = 71333.336
$ws.Range("L63").Value = 71333.336
$ws.Range("N63").Value = -72581.336

$ws.Range("H66").Value = 71333.336
$ws.Range("J66").Value = 71333.336
$ws.Range("L66").Value = 214000.008
$ws.Range("N66").Value = -220240.008

$ws.Range("H81").Value = 4484.2
$ws.Range("I81").Value = 4484.2
$ws.Range("K81").Value = 8968.4
$ws.Range("M81").Value = -7907.4

$ws.Range("H84").Value = 4484.2
$ws.Range("I84").Value = 4484.2
$ws.Range("K84").Value = 44842
$ws.Range("M84").Value = -39538

$ws.Range("H107").Value = 1050.8
$ws.Range("I107").Value = 1160
$ws.Range("J107").Value = 614
$ws.Range("K107").Value = 3480
$ws.Range("L107").Value = 1842
$ws.Range("M107").Value = -1560
$ws.Range("N107").Value = -5682

$ws.Range("H122").Value = 3600.1482
$ws.Range("I122").Value = 2966.25
$ws.Range("K122").Value = 8898.75
$ws.Range("M122").Value = -6448.75

$ws.Range("H136").Value = 15217.286
$ws.Range("I136").Value = 19448
$ws.Range("J136").Value = 1679
$ws.Range("K136").Value = 58344
$ws.Range("L136").Value = 5037
$ws.Range("M136").Value = -55794
$ws.Range("N136").Value = -10137
